$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Multivalued" column header in L4, bold black Calibri 11 (matches other header cells' font,
# but with an explicit color so it becomes its own font entry).
$ws.Range("L4").Value = "Multivalued"
$ws.Range("L4").Font.Bold = $true
$ws.Range("L4").Font.Name = "Calibri"
$ws.Range("L4").Font.Size = 11
$ws.Range("L4").Font.Color = 0

# New "FALSE" values for the Multivalued column on the three property rows.
# Leading apostrophe forces these to be stored as text (matching the existing
# "FALSE" text cells elsewhere in the sheet) instead of being coerced to booleans.
$ws.Range("L5").Value = "'FALSE"
$ws.Range("L6").Value = "'FALSE"
$ws.Range("L7").Value = "'FALSE"

# Selection ends up on L7 after the edits (matches the authored selection state).
$ws.Range("L7").Select() | Out-Null
